# Actualizacion de horarios - Linea 141 - 251
$wb = $excel.ActiveWorkbook

$nuevaHora = "04:18:53"

# ---------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: $nuevaHora"
$ws1.Range("A3").Value = "Total filas: 10"

$ws1.Range("A6").Value = $nuevaHora
$ws1.Range("B6").Value = "04:45"
$ws1.Range("C6").Value = "215A_EL PATO"
$ws1.Range("D6").Value = 27
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $nuevaHora
$ws1.Range("B7").Value = "04:53"
$ws1.Range("C7").Value = "11_ETCHEVERRY"
$ws1.Range("D7").Value = 35
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = $nuevaHora
$ws1.Range("B8").Value = "05:16"
$ws1.Range("C8").Value = "17_ROMERO"
$ws1.Range("D8").Value = 58
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = $nuevaHora
$ws1.Range("B9").Value = "05:21"
$ws1.Range("C9").Value = "23_HERNANDEZ"
$ws1.Range("D9").Value = 63
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = $nuevaHora
$ws1.Range("B10").Value = "05:34"
$ws1.Range("C10").Value = "215B_EL PATO"
$ws1.Range("D10").Value = 76
$ws1.Range("E10").Value = "LP1912"

$ws1.Range("A11").Value = $nuevaHora
$ws1.Range("B11").Value = "05:37"
$ws1.Range("C11").Value = "14_ABASTO"
$ws1.Range("D11").Value = 79
$ws1.Range("E11").Value = "LP1912"

$ws1.Range("A12").Value = $nuevaHora
$ws1.Range("B12").Value = "05:46"
$ws1.Range("C12").Value = "15_ABASTO"
$ws1.Range("D12").Value = 88
$ws1.Range("E12").Value = "LP1912"

$ws1.Range("A13").Value = $nuevaHora
$ws1.Range("B13").Value = "06:07"
$ws1.Range("C13").Value = "16_SANTA ANA"
$ws1.Range("D13").Value = 109
$ws1.Range("E13").Value = "LP1912"

$ws1.Range("A14").Value = $nuevaHora
$ws1.Range("B14").Value = "06:11"
$ws1.Range("C14").Value = "215A_EL PATO"
$ws1.Range("D14").Value = 113
$ws1.Range("E14").Value = "LP1912"

$ws1.Range("A15").Value = $nuevaHora
$ws1.Range("B15").Value = "06:13"
$ws1.Range("C15").Value = "225_HARAS DEL SUR"
$ws1.Range("D15").Value = 115
$ws1.Range("E15").Value = "LP1912"

# ---------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: $nuevaHora"
$ws2.Range("A3").Value = "Total filas: 3"

$ws2.Range("A6").Value = $nuevaHora
$ws2.Range("B6").Value = "04:45"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 27
$ws2.Range("E6").Value = "LP1912"

$ws2.Range("A7").Value = $nuevaHora
$ws2.Range("B7").Value = "05:34"
$ws2.Range("C7").Value = "215B_EL PATO"
$ws2.Range("D7").Value = 76
$ws2.Range("E7").Value = "LP1912"

$ws2.Range("A8").Value = $nuevaHora
$ws2.Range("B8").Value = "06:11"
$ws2.Range("C8").Value = "215A_EL PATO"
$ws2.Range("D8").Value = 113
$ws2.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------
# Hoja 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
$ws3.Range("A3").Value = "Total filas: 2"

# La hoja 3 no tenia encabezado ni filas de datos: copiamos la fila de
# encabezado (valores + formato) desde la hoja 1 para mantener el estilo
$ws1.Range("A5:E5").Copy($ws3.Range("A5:E5"))

$ws3.Range("A6").Value = $nuevaHora
$ws3.Range("B6").Value = "05:43"
$ws3.Range("C6").Value = "215A_LA PLATA"
$ws3.Range("D6").Value = 85
$ws3.Range("E6").Value = "L6173"

$ws3.Range("A7").Value = $nuevaHora
$ws3.Range("B7").Value = "06:08"
$ws3.Range("C7").Value = "215A_LA PLATA"
$ws3.Range("D7").Value = 110
$ws3.Range("E7").Value = "L6173"
